# Apply updated "want to go" counts (column F) and minimum ticket price (column G)
# for sheets "展览" (exhibitions) and "全部类型" (all types).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F5").Value = 1814
$ws1.Range("G5").Value = 54
$ws1.Range("F6").Value = 197
$ws1.Range("F9").Value = 2317
$ws1.Range("F11").Value = 66
$ws1.Range("F14").Value = 500
$ws1.Range("F15").Value = 33
$ws1.Range("F20").Value = 44
$ws1.Range("F24").Value = 77
$ws1.Range("F25").Value = 30
$ws1.Range("F26").Value = 1441
$ws1.Range("F29").Value = 183

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F5").Value = 1814
$ws4.Range("G5").Value = 54
$ws4.Range("F7").Value = 198
$ws4.Range("F10").Value = 2317
$ws4.Range("F12").Value = 66
$ws4.Range("F15").Value = 500
$ws4.Range("F16").Value = 33
$ws4.Range("F21").Value = 44
$ws4.Range("F25").Value = 77
$ws4.Range("F26").Value = 30
$ws4.Range("F27").Value = 1441
$ws4.Range("F30").Value = 184

$wb.Save()
